$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.5853090391950161
$ws.Range("D2").Value = 0.5642960240100428

$ws.Range("C3").Value = -0.9375477372860459
$ws.Range("D3").Value = 0.358656216893062

$ws.Range("C4").Value = 0.2105671085497922
$ws.Range("D4").Value = 0.8351631597310409

$ws.Range("C5").Value = -0.1550878283199024
$ws.Range("D5").Value = 0.8781661814975106

$ws.Range("C6").Value = -0.1200180356651105
$ws.Range("D6").Value = 0.9055581682530507

$ws.Range("C7").Value = 0.5871163470962277
$ws.Range("D7").Value = 0.5631030576436031

$ws.Range("C8").Value = 0.404194577640341
$ws.Range("D8").Value = 0.6899721026778156

$ws.Range("C9").Value = 1.051841271196843
$ws.Range("D9").Value = 0.3042892145689851

$ws.Range("C10").Value = 0.7194869505395031
$ws.Range("D10").Value = 0.4794149507461967

$ws.Range("C11").Value = -0.3632082105602462
$ws.Range("D11").Value = 0.7199159757791866
